$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Global")

# --- Header row (A1:C1) ---
$ws1.Range("A1").Value = "HowMany"
$ws1.Range("B1").Value = "AmountInCart"
$ws1.Range("C1").Value = "Output_Amount_name_out"

# --- Data cells for the new B/C columns (row 2 carries a sample value) ---
$ws1.Range("B2").Value = "CHECKOUT (`$539.98)"
$ws1.Range("C2").Value = "CHECKOUT (`$539.98)"

# --- Borders ---
# Column A (rows 2-5): drop the old right edge, keep top/bottom thin black
# (the cell already carries a thin-black top/bottom from its previous
#  "box" border, so only the right edge needs clearing).
foreach ($r in 2..5) {
    $cell = $ws1.Cells.Item($r, 1)
    $cell.Borders.Item(10).LineStyle = -4142   # xlEdgeRight -> none
}

# Column B (rows 2-5): new thin black top/bottom border (no left/right)
foreach ($r in 2..5) {
    $cell = $ws1.Cells.Item($r, 2)
    $cell.Borders.Item(8).Color = 0
    $cell.Borders.Item(8).LineStyle = 1        # xlEdgeTop
    $cell.Borders.Item(9).Color = 0
    $cell.Borders.Item(9).LineStyle = 1        # xlEdgeBottom
}

# Column C (rows 2-5): thin black top/bottom/right border (closes the row box)
foreach ($r in 2..5) {
    $cell = $ws1.Cells.Item($r, 3)
    $cell.Borders.Item(8).Color = 0
    $cell.Borders.Item(8).LineStyle = 1        # xlEdgeTop
    $cell.Borders.Item(9).Color = 0
    $cell.Borders.Item(9).LineStyle = 1        # xlEdgeBottom
    $cell.Borders.Item(10).Color = 0
    $cell.Borders.Item(10).LineStyle = 1       # xlEdgeRight
}

# --- Column widths: autosize the two new columns to fit their contents ---
$ws1.Columns("B:C").AutoFit()
